$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before the "Late" column (column N), shifting
# "Late" and "Outstanding" one column to the right.
$ws.Range("N1").EntireColumn.Insert()

# Make "Repayment schedule" the active sheet/tab with R7 selected.
$ws.Activate()
$ws.Range("R7").Select()
